# Fill in the "Definition" column (D) for the concept rows on the
# "Concepts" sheet with the same text as the "Display" column (C),
# for rows 2-4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Concepts")

for ($row = 2; $row -le 4; $row++) {
    $display = $ws.Cells.Item($row, 3).Text
    $ws.Cells.Item($row, 4).Value = $display
}
